# Re-positions/resizes the 4 result pictures and repositions 12 of the
# caption/label textboxes on slide 1 (degree_ratio_distributions), per the
# "Adding 5-fold cross-validation on BindingDB data + updated figures" edit.
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points (1 pt =
# 12700 EMU) while the OOXML <a:off>/<a:ext> we need to land on are given in
# EMU. The literal point values below were computed from the target EMU
# values with a tiny half-EMU nudge (sign-matched to the target) so that the
# COM layer's internal float rounding reproduces the exact target EMU
# on save, instead of landing 1 EMU off.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$changes = @{
    7 = @{ Left = 30.876653543307086; Top = 65.10153543307086; Width = 307.12877952755906; Height = 210.32940944881892 }
    9 = @{ Left = 36.62255905511811; Top = 309.7698818897638; Width = 299.4077559055118; Height = 213.00547244094489 }
    11 = @{ Left = 385.56712598425196; Top = 62.44035433070866; Width = 307.09122047244097; Height = 207.55972440944882 }
    13 = @{ Left = 381.2631102362205; Top = 301.60326771653547; Width = 315.699251968504; Height = 226.0406692913386 }
    14 = @{ Left = 171.0478346456693; Top = 53.33027559055118 }
    15 = @{ Left = 118.93594488188977; Top = 257.571062992126 }
    16 = @{ Left = 110.13397637795275; Top = 500.96437007874016 }
    17 = @{ Left = 469.31043307086617; Top = 259.9834251968504 }
    18 = @{ Left = 469.31043307086617; Top = 509.2394881889764 }
    19 = @{ Left = -18.138937007874013; Top = 127.25665354330708 }
    20 = @{ Left = 335.29161417322837; Top = 126.10153543307086 }
    21 = @{ Left = 333.63783464566933; Top = 366.6003543307087 }
    22 = @{ Left = -19.71043307086614; Top = 365.11767716535434 }
    29 = @{ Left = 528.3483858267716; Top = 53.33027559055118 }
    30 = @{ Left = 175.45043307086615; Top = 298.1967322834646 }
    31 = @{ Left = 533.9755511811023; Top = 296.0333464566929 }
}

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($changes.ContainsKey($sh.Id)) {
        $c = $changes[$sh.Id]
        $sh.Left = $c.Left
        $sh.Top = $c.Top
        if ($c.ContainsKey('Width')) { $sh.Width = $c.Width }
        if ($c.ContainsKey('Height')) { $sh.Height = $c.Height }
    }
}
